# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from Excel serial date 45183 (2023-09-14) to 45184 (2023-09-15).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
